# This workbook is refreshed periodically by an automated job that re-pulls
# current market-board data (Universalis) for each Leve sheet and rewrites the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) in place.
# This run updates a handful of rows across all eight crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 771.8
$ws.Range("I6").Value = 771.8
$ws.Range("K6").Value = 2315.4
$ws.Range("M6").Value = -2203.4
$ws.Range("H58").Value = 1371.875
$ws.Range("I58").Value = 38
$ws.Range("J58").Value = 1816.5
$ws.Range("K58").Value = 114
$ws.Range("L58").Value = 5449.5
$ws.Range("M58").Value = 36
$ws.Range("N58").Value = -5749.5
$ws.Range("H94").Value = 876.3333
$ws.Range("I94").Value = 876.3333
$ws.Range("K94").Value = 876.3333
$ws.Range("M94").Value = -425.3333
$ws.Range("H100").Value = 1873.9166
$ws.Range("I100").Value = 1284
$ws.Range("J100").Value = 2699.8
$ws.Range("K100").Value = 1284
$ws.Range("L100").Value = 2699.8
$ws.Range("M100").Value = -743
$ws.Range("N100").Value = -3781.8
$ws.Range("H125").Value = 7916.3335
$ws.Range("I125").Value = 4749.25
$ws.Range("J125").Value = 9499.875
$ws.Range("K125").Value = 42743.25
$ws.Range("L125").Value = 85498.875
$ws.Range("M125").Value = -40283.25
$ws.Range("N125").Value = -90418.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4819.7
$ws.Range("I32").Value = 4917.641
$ws.Range("K32").Value = 4917.641
$ws.Range("M32").Value = -4630.641
$ws.Range("H132").Value = 4463.778
$ws.Range("I132").Value = 3220.5862
$ws.Range("K132").Value = 9661.758600000001
$ws.Range("M132").Value = -7131.758600000001
$ws.Range("H133").Value = 132666.33
$ws.Range("J133").Value = 132666.33
$ws.Range("L133").Value = 132666.33
$ws.Range("N133").Value = -137726.33
$ws.Range("H139").Value = 199887.5
$ws.Range("J139").Value = 199887.5
$ws.Range("L139").Value = 199887.5
$ws.Range("N139").Value = -210167.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 668.0769
$ws.Range("J80").Value = 990.2857
$ws.Range("L80").Value = 990.2857
$ws.Range("N80").Value = -2986.2857
$ws.Range("H83").Value = 668.0769
$ws.Range("J83").Value = 990.2857
$ws.Range("L83").Value = 4951.4285
$ws.Range("N83").Value = -14935.4285
$ws.Range("H86").Value = 1529.5555
$ws.Range("I86").Value = 1380.8572
$ws.Range("K86").Value = 1380.8572
$ws.Range("M86").Value = -257.8571999999999
$ws.Range("H89").Value = 1529.5555
$ws.Range("I89").Value = 1380.8572
$ws.Range("K89").Value = 6904.286
$ws.Range("M89").Value = -1288.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5409.3794
$ws.Range("I58").Value = 2263.5293
$ws.Range("K58").Value = 2263.5293
$ws.Range("M58").Value = -2060.5293
$ws.Range("H86").Value = 9141.888999999999
$ws.Range("I86").Value = 8570
$ws.Range("K86").Value = 8570
$ws.Range("M86").Value = -7447
$ws.Range("H89").Value = 9141.888999999999
$ws.Range("I89").Value = 8570
$ws.Range("K89").Value = 42850
$ws.Range("M89").Value = -37234
$ws.Range("H94").Value = 729.9091
$ws.Range("I94").Value = 492.7143
$ws.Range("J94").Value = 1145
$ws.Range("K94").Value = 492.7143
$ws.Range("L94").Value = 1145
$ws.Range("M94").Value = -41.71429999999998
$ws.Range("N94").Value = -2047
$ws.Range("H105").Value = 1803
$ws.Range("I105").Value = 1991.25
$ws.Range("K105").Value = 1991.25
$ws.Range("M105").Value = -244.25
$ws.Range("H107").Value = 533.9375
$ws.Range("I107").Value = 432.77777
$ws.Range("J107").Value = 664
$ws.Range("K107").Value = 432.77777
$ws.Range("L107").Value = 664
$ws.Range("M107").Value = 1487.22223
$ws.Range("N107").Value = -4504
$ws.Range("H122").Value = 5237.1875
$ws.Range("J122").Value = 5715.8335
$ws.Range("L122").Value = 17147.5005
$ws.Range("N122").Value = -22047.5005
$ws.Range("H134").Value = 8151.2
$ws.Range("I134").Value = 6062.4
$ws.Range("J134").Value = 10240
$ws.Range("K134").Value = 18187.2
$ws.Range("L134").Value = 30720
$ws.Range("M134").Value = -15652.2
$ws.Range("N134").Value = -35790
$ws.Range("H136").Value = 5409.3794
$ws.Range("I136").Value = 2263.5293
$ws.Range("K136").Value = 6790.5879
$ws.Range("M136").Value = -4240.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8389.75
$ws.Range("I2").Value = 59
$ws.Range("J2").Value = 12555.125
$ws.Range("K2").Value = 354
$ws.Range("L2").Value = 75330.75
$ws.Range("M2").Value = -241
$ws.Range("N2").Value = -75556.75
$ws.Range("H4").Value = 32147676
$ws.Range("I4").Value = 36698364
$ws.Range("J4").Value = 1203000
$ws.Range("K4").Value = 110095092
$ws.Range("L4").Value = 3609000
$ws.Range("M4").Value = -110094980
$ws.Range("N4").Value = -3609224
$ws.Range("H9").Value = 220
$ws.Range("I9").Value = 225
$ws.Range("K9").Value = 675
$ws.Range("M9").Value = -451
$ws.Range("H16").Value = 61
$ws.Range("J16").Value = 50
$ws.Range("L16").Value = 150
$ws.Range("N16").Value = -496
$ws.Range("H33").Value = 707.0769
$ws.Range("I33").Value = 119.125
$ws.Range("J33").Value = 1647.8
$ws.Range("K33").Value = 714.75
$ws.Range("L33").Value = 9886.799999999999
$ws.Range("M33").Value = -431.75
$ws.Range("N33").Value = -10452.8
$ws.Range("H112").Value = 333336160
$ws.Range("I112").Value = 500001760
$ws.Range("K112").Value = 1500005280
$ws.Range("M112").Value = -1500004172
$ws.Range("H122").Value = 515.75
$ws.Range("J122").Value = 596.3333
$ws.Range("L122").Value = 5366.9997
$ws.Range("N122").Value = -10266.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3170.8948
$ws.Range("I102").Value = 2627.625
$ws.Range("J102").Value = 6068.3335
$ws.Range("K102").Value = 2627.625
$ws.Range("L102").Value = 6068.3335
$ws.Range("M102").Value = -1005.625
$ws.Range("N102").Value = -9312.333500000001
$ws.Range("H126").Value = 2581.6897
$ws.Range("I126").Value = 1801.2941
$ws.Range("J126").Value = 3687.25
$ws.Range("K126").Value = 5403.8823
$ws.Range("L126").Value = 11061.75
$ws.Range("M126").Value = -2933.8823
$ws.Range("N126").Value = -16001.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1591.5
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""
$ws.Range("H74").Value = 42124
$ws.Range("I74").Value = 42124
$ws.Range("K74").Value = 42124
$ws.Range("M74").Value = -41126
$ws.Range("H77").Value = 42124
$ws.Range("I77").Value = 42124
$ws.Range("K77").Value = 126372
$ws.Range("M77").Value = -121380
$ws.Range("H104").Value = 30184.5
$ws.Range("J104").Value = 30184.5
$ws.Range("L104").Value = 30184.5
$ws.Range("N104").Value = -37172.5
$ws.Range("H126").Value = 1591.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5332.6665
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 5332.6665
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H132").Value = 3637.3928
$ws.Range("I132").Value = 2813.3333
$ws.Range("K132").Value = 8439.999899999999
$ws.Range("M132").Value = -5909.999899999999
$ws.Range("H136").Value = 5535.393
$ws.Range("I136").Value = 3449.5
$ws.Range("K136").Value = 10348.5
$ws.Range("M136").Value = -7798.5

Write-Host "Applied 193 updates"
